$wb = $excel.ActiveWorkbook

$sheetNames = @("DisplayValues", "SignificanceValues", "HistDisplayValues", "HistSignificanceValues")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "2;3"
    $ws.Range("C1").Value = "2;4"
}
